# Daily attendance processing - 2025-12-05 22:26:36
# Normalize the "Recorded By" (column G) text so that the "System" token
# is listed first for specific recorder combinations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "admin@admin.com, System"             = "System, admin@admin.com"
    "backup@backdoor.com, system, System" = "backup@backdoor.com, System, system"
}

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $val = $cell.Value2
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value2 = $map[$val]
    }
}
